$wb = $excel.ActiveWorkbook

# --- ALC!row4 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2699.6667
$ws.Range("I4").Value = 49.5
$ws.Range("K4").Value = 49.5
$ws.Range("M4").Value = 64.5

# --- ALC!row29 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("N29").Value = 0

# --- ALC!row45 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 1902.5
$ws.Range("J45").Value = 4610
$ws.Range("L45").Value = 13830
$ws.Range("N45").Value = -14214

# --- ALC!row49 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 3770.111
$ws.Range("I49").Value = 1232.3334
$ws.Range("K49").Value = 3697.0002
$ws.Range("M49").Value = -3561.0002

# --- ALC!row54 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 11388.75
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 11388.75
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("M54").Value = 11388.75
$ws.Range("N54").Value = -12360.75

# --- ALC!row58 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1892.2142
$ws.Range("I58").Value = 610.1111
$ws.Range("J58").Value = 4200
$ws.Range("K58").Value = 1830.3333
$ws.Range("L58").Value = 12600
$ws.Range("M58").Value = -1680.3333
$ws.Range("N58").Value = -12900

# --- ALC!row111 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3631.75
$ws.Range("I111").Value = 2842.3333
$ws.Range("K111").Value = 8526.999899999999
$ws.Range("M111").Value = -5459.999899999999

# --- ALC!row138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3817.9324
$ws.Range("I138").Value = 1643.6086
$ws.Range("J138").Value = 4798.51
$ws.Range("K138").Value = 4930.825800000001
$ws.Range("L138").Value = 14395.53
$ws.Range("M138").Value = 209.1741999999995
$ws.Range("N138").Value = -24675.53

# --- ARM!row32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6216.3237
$ws.Range("I32").Value = 5987.4194
$ws.Range("K32").Value = 5987.4194
$ws.Range("M32").Value = -5700.4194

# --- ARM!row102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 11833.728
$ws.Range("I102").Value = 13475.223
$ws.Range("K102").Value = 13475.223
$ws.Range("M102").Value = -11853.223

# --- ARM!row132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4465.1196
$ws.Range("J132").Value = 6997.5
$ws.Range("L132").Value = 20992.5
$ws.Range("N132").Value = -26052.5

# --- BSM!row7 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2506811.2
$ws.Range("I7").Value = 6124
$ws.Range("J7").Value = 5007498.5
$ws.Range("K7").Value = 6124
$ws.Range("L7").Value = 5007498.5
$ws.Range("M7").Value = -6011
$ws.Range("N7").Value = -5007724.5

# --- BSM!row97 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 42000
$ws.Range("I97").Value = 12500
$ws.Range("K97").Value = 12500
$ws.Range("M97").Value = -11509

# --- BSM!row134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9441.457
$ws.Range("I134").Value = 9831.733
$ws.Range("J134").Value = 7099.8
$ws.Range("K134").Value = 29495.199
$ws.Range("L134").Value = 21299.4
$ws.Range("M134").Value = -26960.199
$ws.Range("N134").Value = -26369.4

# --- BSM!row140 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 78623.84
$ws.Range("J140").Value = 77274.28999999999
$ws.Range("L140").Value = 77274.28999999999
$ws.Range("N140").Value = -87634.28999999999

# --- CRP!row16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1191.7646
$ws.Range("I16").Value = 1188.3334
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1188.3334
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -901.3334
$ws.Range("N16").Value = -1774

# --- CRP!row58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4482.8335
$ws.Range("I58").Value = 4279.4
$ws.Range("J58").Value = 5500
$ws.Range("K58").Value = 4279.4
$ws.Range("L58").Value = 5500
$ws.Range("M58").Value = -4076.4
$ws.Range("N58").Value = -5906

# --- CRP!row95 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 104034900
$ws.Range("J95").Value = 104034900
$ws.Range("L95").Value = 104034900
$ws.Range("N95").Value = -104040392

# --- CRP!row113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1191.7646
$ws.Range("I113").Value = 1188.3334
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1188.3334
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 981.6666
$ws.Range("N113").Value = -5540

# --- CRP!row132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10002.32
$ws.Range("I132").Value = 1498.4359
$ws.Range("J132").Value = 40152.453
$ws.Range("K132").Value = 4495.307699999999
$ws.Range("L132").Value = 120457.359
$ws.Range("M132").Value = -1965.307699999999
$ws.Range("N132").Value = -125517.359

# --- CRP!row136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4482.8335
$ws.Range("I136").Value = 4279.4
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 12838.2
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -10288.2
$ws.Range("N136").Value = -21600

# --- CUL!row17 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2096.8572
$ws.Range("I17").Value = 1132
$ws.Range("J17").Value = 2820.5
$ws.Range("K17").Value = 3396
$ws.Range("L17").Value = 8461.5
$ws.Range("M17").Value = -3227
$ws.Range("N17").Value = -8799.5

# --- CUL!row29 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 184.54546
$ws.Range("I29").Value = 210
$ws.Range("K29").Value = 630
$ws.Range("M29").Value = -353

# --- CUL!row36 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 125
$ws.Range("I36").Value = 125
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 375
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -206

# --- GSM!row2 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1099.8125
$ws.Range("I2").Value = 1378.9166
$ws.Range("K2").Value = 1378.9166
$ws.Range("M2").Value = -1265.9166

# --- GSM!row70 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13267.667
$ws.Range("I70").Value = 12849.5
$ws.Range("K70").Value = 12849.5
$ws.Range("M70").Value = -12579.5

# --- GSM!row73 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13267.667
$ws.Range("I73").Value = 12849.5
$ws.Range("K73").Value = 12849.5
$ws.Range("M73").Value = -11913.5

# --- GSM!row102 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5451.3657
$ws.Range("J102").Value = 2445.182
$ws.Range("L102").Value = 2445.182
$ws.Range("N102").Value = -5689.182

# --- LTW!row22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8043
$ws.Range("J22").Value = 5175.25
$ws.Range("L22").Value = 5175.25
$ws.Range("N22").Value = -5765.25

# --- LTW!row27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8043
$ws.Range("J27").Value = 5175.25
$ws.Range("L27").Value = 5175.25
$ws.Range("N27").Value = -5389.25

# --- LTW!row93 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8176.8
$ws.Range("I93").Value = 8176.8
$ws.Range("K93").Value = 8176.8
$ws.Range("M93").Value = -6928.8

# --- LTW!row122 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6583.9414
$ws.Range("I122").Value = 4866.385
$ws.Range("J122").Value = 12166
$ws.Range("K122").Value = 14599.155
$ws.Range("L122").Value = 36498
$ws.Range("M122").Value = -12149.155
$ws.Range("N122").Value = -41398

# --- WVR!row141 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 105000
$ws.Range("J141").Value = 105000
$ws.Range("L141").Value = 105000
$ws.Range("N141").Value = -115360
